$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (20) of forecast data, matching the pattern of row 19.
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 1.049317648994741
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 0.72625340902297

# Carry over the date cell's number/style formatting from the row above (A19 -> A20).
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null

$ws.Range("A1").Select() | Out-Null
